$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "07-sep" column (BL) following the existing "06-sep" column (BK)
# Header cell BL1, formatted like BK1 (text style)
$ws.Range("BK1").Copy()
$ws.Range("BL1").PasteSpecial(-4122)
$ws.Range("BL1").Value = "07-sep"

# Data cells BL2:BL11, formatted like BK2:BK11 (centered number style)
$ws.Range("BK2:BK11").Copy()
$ws.Range("BL2:BL11").PasteSpecial(-4122)

$ws.Range("BL2").Value = 17
$ws.Range("BL3").Value = 14
$ws.Range("BL4").Value = 11
$ws.Range("BL5").Value = 13
$ws.Range("BL6").Value = 8
$ws.Range("BL7").Value = 18
$ws.Range("BL8").Value = 24
$ws.Range("BL9").Value = 14
$ws.Range("BL10").Value = 15
$ws.Range("BL11").Value = 18

# Update the active selection to match the newly extended range
$ws.Range("BL12").Select()
